$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data (refreshed figures from GitHub Actions run)

$ws.Range("D2").Value = "28.505.25"
$ws.Range("E2").Value = "  +0.14%  "

$ws.Range("D3").Value = "1.870.52"
$ws.Range("E3").Value = "  -0.11%  "

$ws.Range("D4").Formula = "=""1.007"""
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D4").PasteSpecial(-4163) | Out-Null
$ws.Range("E4").Value = "  -1.29%  "

$ws.Range("D5").Formula = "=""315.31"""
$ws.Range("D5").Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4163) | Out-Null
$ws.Range("E5").Value = "  -0.52%  "

$ws.Range("D6").Formula = "=""1.006"""
$ws.Range("D6").Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4163) | Out-Null
$ws.Range("E6").Value = "  -1.38%  "

$ws.Range("D7").Formula = "=""0.5079"""
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D7").PasteSpecial(-4163) | Out-Null
$ws.Range("E7").Value = "  -1.05%  "

$ws.Range("D8").Formula = "=""0.3904"""
$ws.Range("D8").Copy() | Out-Null
$ws.Range("D8").PasteSpecial(-4163) | Out-Null
$ws.Range("E8").Value = "  -0.91%  "

$ws.Range("D9").Formula = "=""0.08369"""
$ws.Range("D9").Copy() | Out-Null
$ws.Range("D9").PasteSpecial(-4163) | Out-Null
$ws.Range("E9").Value = "  +0.61%  "

$ws.Range("E10").Value = "  -0.70%  "

$ws.Range("D11").Formula = "=""41.77"""
$ws.Range("D11").Copy() | Out-Null
$ws.Range("D11").PasteSpecial(-4163) | Out-Null
$ws.Range("E11").Value = "  -0.49%  "

$ws.Range("D12").Formula = "=""6.222"""
$ws.Range("D12").Copy() | Out-Null
$ws.Range("D12").PasteSpecial(-4163) | Out-Null
$ws.Range("E12").Value = "  -0.17%  "

$ws.Range("D13").Value = "1.873.04"
$ws.Range("E13").Value = "  +0.71%  "

$ws.Range("D14").Formula = "=""20.45"""
$ws.Range("D14").Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4163) | Out-Null
$ws.Range("E14").Value = "  -0.06%  "

$ws.Range("D15").Formula = "=""7.283"""
$ws.Range("D15").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4163) | Out-Null
$ws.Range("E15").Value = "  +0.84%  "

$ws.Range("E16").Value = "  -1.35%  "

$ws.Range("D17").Formula = "=""0.00001101"""
$ws.Range("D17").Copy() | Out-Null
$ws.Range("D17").PasteSpecial(-4163) | Out-Null
$ws.Range("E17").Value = "  -0.55%  "

$ws.Range("D18").Formula = "=""91.15"""
$ws.Range("D18").Copy() | Out-Null
$ws.Range("D18").PasteSpecial(-4163) | Out-Null
$ws.Range("E18").Value = "  -0.09%  "

$ws.Range("D19").Formula = "=""0.06728"""
$ws.Range("D19").Copy() | Out-Null
$ws.Range("D19").PasteSpecial(-4163) | Out-Null
$ws.Range("E19").Value = "  -0.22%  "

$ws.Range("E20").Value = "  +0.61%  "

$ws.Range("E21").Value = "  -1.39%  "

$ws.Range("D22").Formula = "=""5.928"""
$ws.Range("D22").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4163) | Out-Null
$ws.Range("E22").Value = "  -0.51%  "

$ws.Range("D23").Value = "28.527.20"
$ws.Range("E23").Value = "  +0.03%  "

$ws.Range("D24").Formula = "=""11.11"""
$ws.Range("D24").Copy() | Out-Null
$ws.Range("D24").PasteSpecial(-4163) | Out-Null
$ws.Range("E24").Value = "  -0.24%  "

$ws.Range("D25").Formula = "=""2.206"""
$ws.Range("D25").Copy() | Out-Null
$ws.Range("D25").PasteSpecial(-4163) | Out-Null
$ws.Range("E25").Value = "  -2.47%  "

$ws.Range("D26").Value = "2.084.36"
$ws.Range("E26").Value = "  +0.50%  "

$ws.Range("D27").Formula = "=""160.29"""
$ws.Range("D27").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4163) | Out-Null
$ws.Range("E27").Value = "  -0.84%  "

$ws.Range("E28").Value = "  -0.68%  "

$ws.Range("D29").Formula = "=""2.429"""
$ws.Range("D29").Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4163) | Out-Null
$ws.Range("E29").Value = "  +2.68%  "

$ws.Range("D30").Formula = "=""127.27"""
$ws.Range("D30").Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4163) | Out-Null
$ws.Range("E30").Value = "  +0.18%  "

$ws.Range("D31").Formula = "=""0.1039"""
$ws.Range("D31").Copy() | Out-Null
$ws.Range("D31").PasteSpecial(-4163) | Out-Null
$ws.Range("E31").Value = "  -1.46%  "

$ws.Range("E32").Value = "  +0.62%  "

$ws.Range("D33").Formula = "=""5.745"""
$ws.Range("D33").Copy() | Out-Null
$ws.Range("D33").PasteSpecial(-4163) | Out-Null
$ws.Range("E33").Value = "  -1.09%  "

$ws.Range("D34").Formula = "=""3.626"""
$ws.Range("D34").Copy() | Out-Null
$ws.Range("D34").PasteSpecial(-4163) | Out-Null
$ws.Range("E34").Value = "  -0.58%  "

$ws.Range("D35").Formula = "=""0.02457"""
$ws.Range("D35").Copy() | Out-Null
$ws.Range("D35").PasteSpecial(-4163) | Out-Null
$ws.Range("E35").Value = "  +0.77%  "

$ws.Range("D36").Formula = "=""0.06584"""
$ws.Range("D36").Copy() | Out-Null
$ws.Range("D36").PasteSpecial(-4163) | Out-Null
$ws.Range("E36").Value = "  +1.25%  "

$ws.Range("D37").Formula = "=""8.926"""
$ws.Range("D37").Copy() | Out-Null
$ws.Range("D37").PasteSpecial(-4163) | Out-Null
$ws.Range("E37").Value = "  -2.51%  "

$ws.Range("D38").Formula = "=""0.2166"""
$ws.Range("D38").Copy() | Out-Null
$ws.Range("D38").PasteSpecial(-4163) | Out-Null
$ws.Range("E38").Value = "  -0.72%  "

$ws.Range("D39").Formula = "=""5.042"""
$ws.Range("D39").Copy() | Out-Null
$ws.Range("D39").PasteSpecial(-4163) | Out-Null
$ws.Range("E39").Value = "  +1.03%  "

$ws.Range("D40").Formula = "=""1.181"""
$ws.Range("D40").Copy() | Out-Null
$ws.Range("D40").PasteSpecial(-4163) | Out-Null
$ws.Range("E40").Value = "  -0.24%  "

$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Formula = "=""0.6370"""
$ws.Range("D41").Copy() | Out-Null
$ws.Range("D41").PasteSpecial(-4163) | Out-Null
$ws.Range("E41").Value = "  -1.05%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Formula = "=""1.231"""
$ws.Range("D42").Copy() | Out-Null
$ws.Range("D42").PasteSpecial(-4163) | Out-Null
$ws.Range("E42").Value = "  -1.40%  "

$ws.Range("D43").Formula = "=""11.10"""
$ws.Range("D43").Copy() | Out-Null
$ws.Range("D43").PasteSpecial(-4163) | Out-Null
$ws.Range("E43").Value = "  -0.46%  "

$ws.Range("D44").Formula = "=""1.006"""
$ws.Range("D44").Copy() | Out-Null
$ws.Range("D44").PasteSpecial(-4163) | Out-Null
$ws.Range("E44").Value = "  -1.25%  "

$ws.Range("D45").Formula = "=""0.6014"""
$ws.Range("D45").Copy() | Out-Null
$ws.Range("D45").PasteSpecial(-4163) | Out-Null
$ws.Range("E45").Value = "  -0.23%  "

$ws.Range("D46").Formula = "=""13.00"""
$ws.Range("D46").Copy() | Out-Null
$ws.Range("D46").PasteSpecial(-4163) | Out-Null
$ws.Range("E46").Value = "  -0.20%  "

$ws.Range("D47").Formula = "=""3.683"""
$ws.Range("D47").Copy() | Out-Null
$ws.Range("D47").PasteSpecial(-4163) | Out-Null
$ws.Range("E47").Value = "  -0.44%  "

$ws.Range("D48").Formula = "=""2.006"""
$ws.Range("D48").Copy() | Out-Null
$ws.Range("D48").PasteSpecial(-4163) | Out-Null
$ws.Range("E48").Value = "  +0.88%  "

$ws.Range("D49").Formula = "=""1.217"""
$ws.Range("D49").Copy() | Out-Null
$ws.Range("D49").PasteSpecial(-4163) | Out-Null
$ws.Range("E49").Value = "  +0.75%  "

$ws.Range("D50").Formula = "=""122.42"""
$ws.Range("D50").Copy() | Out-Null
$ws.Range("D50").PasteSpecial(-4163) | Out-Null
$ws.Range("E50").Value = "  +0.35%  "

$ws.Range("D51").Formula = "=""0.06809"""
$ws.Range("D51").Copy() | Out-Null
$ws.Range("D51").PasteSpecial(-4163) | Out-Null
$ws.Range("E51").Value = "  -0.71%  "

$excel.CutCopyMode = $false